# Generate Report for Handback
# Update timestamp cells to reflect the latest handback/generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2, column G)
$wsOverview.Range("G2").Value = "2016-10-19 17:27:12"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-10-19 17:27:01"
$wsZhCn.Range("K2").Value = "2016-10-19 17:27:39"

# de-de sheet: "Correspond Handback DateTime" (K2)
$wsDeDe.Range("K2").Value = "2016-10-19 17:27:57"
